# Move Question 96 from the Audit Report group to the Year-End Report
# group (data-entry fix), and nudge the UI state (selection / tab ratio)
# the way the original author's Excel session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuestionsGroups")

# --- Year-End Report row (row 15): questions were "84-95, 139" -------------
# Question 96 now belongs here, so the range grows to include it.
$ws.Cells.Item(15, 3).Value = "84-96, 139"

# --- Audit Report row (row 16): questions were "96-102, 140" ---------------
# Question 96 moves out, so the range now starts at 97.
$ws.Cells.Item(16, 3).Value = "97-102, 140"

# --- Blank spacer rows gain an explicitly formatted cell in column A -------
# (matches the formatting already present on the other blank spacer rows,
# e.g. row 2 above the first group and row 18 above the last group).
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("A18").Copy()
$ws.Range("A17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Window / selection state ----------------------------------------------
# The author's session ended with the active cell on C17 and a slightly
# wider sheet-tab strip.
$win = $excel.ActiveWindow
$win.TabRatio = 0.992

$ws.Range("C17").Select()
